# Fruta / hortaliza, semanal
# Insert a new daily record row at row 272 (pushing existing rows 272-320 down to 273-321)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 272; this shifts rows 272-320 -> 273-321
# and Excel extends the used range / dimension automatically (A1:T320 -> A1:T321).
$ws.Rows.Item(272).Insert()

# Populate the constant (A-K) columns identically to the rest of the dataset.
$ws.Range("A272").Value = 4
$ws.Range("B272").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C272").Value = "Los Lagos"
$ws.Range("D272").Value = 44522
$ws.Range("E272").Value = 10
$ws.Range("F272").Value = "Fruta"
$ws.Range("G272").Value = 100102
$ws.Range("H272").Value = "Cítricos"
$ws.Range("I272").Value = 100102003
$ws.Range("J272").Value = "Limón"
$ws.Range("K272").Value = "Sin especificar"

# Populate the new observation's specific values.
$ws.Range("L272").Value = "1a amarillo"
$ws.Range("M272").Value = 400
$ws.Range("N272").Value = 12000
$ws.Range("O272").Value = 12500
$ws.Range("P272").Value = 12250
$ws.Range("Q272").Value = "$/malla 18 kilos"
$ws.Range("R272").Value = "Provincia de Curicó"
$ws.Range("S272").Value = 681
$ws.Range("T272").Value = 18
